# Add a new slide by duplicating the existing Title Slide (same layout:
# ctrTitle + subTitle placeholders) and updating its title text, matching
# the "Add files via upload" commit that appended a second slide.

$p = $ppt.ActivePresentation

$sourceSlide = $p.Slides.Item(1)
$newSlideRange = $sourceSlide.Duplicate()
$newSlide = $newSlideRange.Item(1)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "test2"
